# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted at row 74 (pushing the existing
# rows 74-148 down to 75-149) on the single worksheet of the workbook.
# dimension grows from A1:R148 to A1:R149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 74; everything below (old rows 74-148)
# shifts down to 75-149, carrying its own values/styles with it.
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with this week's record.
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 44897
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = 100112021
$ws.Range("G74").Value = "Ají"
$ws.Range("H74").Value = "Americana (o)"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 50
$ws.Range("K74").Value = 40000
$ws.Range("L74").Value = 42000
$ws.Range("M74").Value = 41200
$ws.Range("N74").Value = "$/caja 25 kilos"
$ws.Range("O74").Value = "Provincia de Limarí"
$ws.Range("P74").Value = 1648
$ws.Range("Q74").Value = 25
$ws.Range("R74").Value = "Hortaliza"
